$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.764.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.033"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.033"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4427"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3813"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07468"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.77"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.52"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.554"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.787"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07202"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.038"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009141"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.032"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.778.91"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.324"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.32"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.093.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.020"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.64"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.376"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.989"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09049"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.235"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7813"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.028"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.596"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.035"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.145"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01994"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05369"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.881"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5226"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1695"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.928"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.743"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "111.62"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06664"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.037"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.717"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4736"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.926"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.52%  "
